# Update the "ランサーズ" worksheet with the freshly scraped listing data
# (commit message: "Append: 2025-09-12 06:27 JST").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newDate = "2025-09-12 06:27:00"

# --- Row 2 : only the timestamp changes, everything else stays the same ---
$ws.Range("A2").Value = $newDate

# --- Row 3 ---
$ws.Range("A3").Value = $newDate
$ws.Range("B3").Value = "【全世界1億DL】ペイントアプリ『MediBang Paint』のAndroidエンジニア"
$ws.Range("D3").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5392277"
$ws.Range("G3").Value = 345
$ws.Range("H3").Value = "🔥AI,Ai ◇アプリ"

# --- Row 4 ---
$ws.Range("A4").Value = $newDate
$ws.Range("B4").Value = "GPUサーバー導入による節税シミュレーションツール開発"
$ws.Range("D4").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5392249"
$ws.Range("G4").Value = 123
$ws.Range("H4").Value = "◆ツール,開発"

# --- Row 5 ---
$ws.Range("A5").Value = $newDate
$ws.Range("B5").Value = "【急募】EA自動化システム構築の専門家を探しています!"
$ws.Range("D5").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5392078"
$ws.Range("G5").Value = 95
$ws.Range("H5").Value = "◆自動化"

# --- Row 6 ---
$ws.Range("A6").Value = $newDate
$ws.Range("B6").Value = "サブスク型学習サイトの開発"
$ws.Range("D6").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5365024"
$ws.Range("G6").Value = 93
$ws.Range("H6").Value = "◆開発 ◇サイト"

# --- Row 7 ---
$ws.Range("A7").Value = $newDate
$ws.Range("B7").Value = "IB報酬を得るための高性能EA開発依頼"
$ws.Range("D7").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5392235"
$ws.Range("G7").Value = 68
$ws.Range("H7").Value = "◆開発"

# --- Row 8 ---
$ws.Range("A8").Value = $newDate
$ws.Range("B8").Value = "【急募】在庫管理システムの構築!(その後手配管理システムも依頼予定)"
$ws.Range("D8").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5392325"
$ws.Range("G8").Value = 53
$ws.Range("H8").Value = "◇管理"

# --- Row 9 (no skill-tag column anymore) ---
$ws.Range("A9").Value = $newDate
$ws.Range("B9").Value = "Vue.jsを使用した「既存ページ修正」+「追加実装」(ピクセルパーフェクト実装)"
$ws.Range("D9").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5392236"
$ws.Range("G9").Value = 18
$ws.Range("H9").Value = ""

# --- Row 10 (no skill-tag column anymore) ---
$ws.Range("A10").Value = $newDate
$ws.Range("B10").Value = "【急募】A1活用 画像加工とCSV作成のプロを探しています!"
$ws.Range("D10").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5392360"
$ws.Range("G10").Value = 13
$ws.Range("H10").Value = ""

# --- Row 11 (no skill-tag column anymore) ---
$ws.Range("A11").Value = $newDate
$ws.Range("B11").Value = "【急募】エクセルVBAからXLLアドイン作成の依頼"
$ws.Range("D11").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5392307"
$ws.Range("G11").Value = 10
$ws.Range("H11").Value = ""

# --- Drop the old rows 12-18 entirely, shrinking the sheet to A1:H11 ---
$ws.Rows("12:18").Delete()

# --- Column width tweaks (B: 49 -> 47, D: 30 -> 32) ---
# NOTE: the ColumnWidth COM property and the OOXML stored width are not on a
# 1:1 scale in this engine, so we compensate by a small empirically
# determined offset so that the persisted <col width="..."/> lands exactly
# on the intended integer value.
$ws.Columns.Item(2).ColumnWidth = 47 - 0.9
$ws.Columns.Item(4).ColumnWidth = 32 - 0.9

# --- Rebuild the hyperlinks so only F2:F11 keep a live link (old F12:F18
# links must disappear, and F3:F11 must point at the new work URLs) ---
$firstHyperlink = $ws.Hyperlinks.Item(1)
$firstHyperlink.Range.Hyperlinks.Delete()

$links = @(
    @{ Cell = "F2";  Url = "https://www.lancers.jp/work/detail/5392099" },
    @{ Cell = "F3";  Url = "https://www.lancers.jp/work/detail/5392277" },
    @{ Cell = "F4";  Url = "https://www.lancers.jp/work/detail/5392249" },
    @{ Cell = "F5";  Url = "https://www.lancers.jp/work/detail/5392078" },
    @{ Cell = "F6";  Url = "https://www.lancers.jp/work/detail/5365024" },
    @{ Cell = "F7";  Url = "https://www.lancers.jp/work/detail/5392235" },
    @{ Cell = "F8";  Url = "https://www.lancers.jp/work/detail/5392325" },
    @{ Cell = "F9";  Url = "https://www.lancers.jp/work/detail/5392236" },
    @{ Cell = "F10"; Url = "https://www.lancers.jp/work/detail/5392360" },
    @{ Cell = "F11"; Url = "https://www.lancers.jp/work/detail/5392307" }
)

foreach ($link in $links) {
    $rng = $ws.Range($link.Cell)
    $ws.Hyperlinks.Add($rng, $link.Url)
    $rng.Style = "Hyperlink"
}
